# Condition Two-Stage Letter.docx - apply 2021 NPPF/London Plan revisions
# + a handful of wording/font tweaks, as described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

function Set-RunFont($findText) {
    # Locate the (unique) text and tag it with the Source Sans Pro Light
    # ascii/hAnsi font, keeping whatever w:cs it already has.
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if ($found) {
        $rng.Font.NameAscii = "Source Sans Pro Light"
        $rng.Font.NameOther = "Source Sans Pro Light"
    }
}

# --- Wording updates -------------------------------------------------

Replace-Text "NATIONAL PLANNING POLICY FRAMEWORK 2019" "NATIONAL PLANNING POLICY FRAMEWORK 2021"

Replace-Text "NPPF Section 16 and the Draft London Plan (2017 Policy HC1)" "NPPF Section 16 and the London Plan (2021 Policy HC1)"

Replace-Text "NPPF paragraph 189 says applicants" "NPPF paragraph 194 says applicants"

Replace-Text "NPPF paragraphs 185 and 192 and Draft London Plan Policy HC1" "NPPF paragraphs 190 and 197 and London Plan Policy HC1"

Replace-Text "If you grant planning consent, paragraph 199 of the NPPF" "If you grant planning consent, paragraph 205 of the NPPF"

Replace-Text "this pre-commencement condition please let us know" "this pre-commencement condition, please let us know"

Replace-Text "being imposed the application should be refused as it would not comply with NPPF paragraph 199." "being imposed the application should be refused as it would not comply with NPPF paragraph 205."

Replace-Text "Archaeology Advisor" "Archaeology Adviser"

# --- Font-only tweaks (Arial-only runs become Source Sans Pro Light) --

Set-RunFont "Reference"
Set-RunFont "Primary Reference Number"
Set-RunFont "Contact: <Casework Officer>"
Set-RunFont "Casework Officer Number>"
Set-RunFont "Casework Officer Email"
Set-RunFont "Completion Date>"
Set-RunFont "Contact Name>"
Set-RunFont "Log Date>"
